$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 1543.3125
$ws.Range("I19").Value = 2027.5
$ws.Range("J19").Value = 736.3333
$ws.Range("K19").Value = 2027.5
$ws.Range("L19").Value = 736.3333
$ws.Range("M19").Value = -1852.5
$ws.Range("N19").Value = -1086.3333

$ws.Range("H28").Value = 1057.7142
$ws.Range("I28").Value = 981.8
$ws.Range("J28").Value = 1247.5
$ws.Range("K28").Value = 981.8
$ws.Range("L28").Value = 1247.5
$ws.Range("M28").Value = -496.8
$ws.Range("N28").Value = -2217.5

$ws.Range("H63").Value = 66562.60000000001
$ws.Range("I63").Value = 46000
$ws.Range("K63").Value = 46000
$ws.Range("M63").Value = -45376

$ws.Range("H66").Value = 66562.60000000001
$ws.Range("I66").Value = 46000
$ws.Range("K66").Value = 138000
$ws.Range("M66").Value = -134880

$ws.Range("H137").Value = 32260296
$ws.Range("I137").Value = 52633370
$ws.Range("K137").Value = 157900110
$ws.Range("M137").Value = -157897560

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6331522
$ws.Range("I32").Value = 7043878.5
$ws.Range("K32").Value = 7043878.5
$ws.Range("M32").Value = -7043591.5

$ws.Range("H122").Value = 0
$ws.Range("I122").Value = 0
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 0
$ws.Range("L122").Value = 0
$ws.Range("M122").ClearContents()
$ws.Range("N122").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H9").Value = 100000
$ws.Range("J9").Value = 100000
$ws.Range("L9").Value = 100000
$ws.Range("N9").Value = -100336

$ws.Range("H16").Value = 3499.5
$ws.Range("J16").Value = 4999
$ws.Range("L16").Value = 4999
$ws.Range("N16").Value = -5573

$ws.Range("H31").Value = 78089.14999999999
$ws.Range("I31").Value = 93982.66
$ws.Range("K31").Value = 93982.66
$ws.Range("M31").Value = -93687.66

$ws.Range("H34").Value = 78089.14999999999
$ws.Range("I34").Value = 93982.66
$ws.Range("K34").Value = 93982.66
$ws.Range("M34").Value = -93780.66

$ws.Range("H39").Value = 8864
$ws.Range("J39").Value = 9999.25
$ws.Range("L39").Value = 9999.25
$ws.Range("N39").Value = -10781.25

$ws.Range("H41").Value = 22649.572
$ws.Range("I41").Value = 2208
$ws.Range("J41").Value = 34006
$ws.Range("K41").Value = 2208
$ws.Range("L41").Value = 34006
$ws.Range("M41").Value = -1780
$ws.Range("N41").Value = -34862

$ws.Range("H49").Value = 8864
$ws.Range("J49").Value = 9999.25
$ws.Range("L49").Value = 9999.25
$ws.Range("N49").Value = -10363.25

$ws.Range("H59").Value = 97399.39999999999
$ws.Range("J59").Value = 97399.39999999999
$ws.Range("L59").Value = 97399.39999999999
$ws.Range("N59").Value = -99689.39999999999

$ws.Range("H99").Value = 4277.5713
$ws.Range("I99").Value = 3610.75
$ws.Range("J99").Value = 5166.6665
$ws.Range("K99").Value = 3610.75
$ws.Range("L99").Value = 5166.6665
$ws.Range("M99").Value = -2112.75
$ws.Range("N99").Value = -8162.6665

$ws.Range("H113").Value = 3499.5
$ws.Range("J113").Value = 4999
$ws.Range("L113").Value = 4999
$ws.Range("N113").Value = -9339

$ws.Range("H126").Value = 4277.5713
$ws.Range("I126").Value = 3610.75
$ws.Range("J126").Value = 5166.6665
$ws.Range("K126").Value = 10832.25
$ws.Range("L126").Value = 15499.9995
$ws.Range("M126").Value = -8362.25
$ws.Range("N126").Value = -20439.9995

$ws.Range("H132").Value = 45249188
$ws.Range("I132").Value = 43023184
$ws.Range("K132").Value = 129069552
$ws.Range("M132").Value = -129067022

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 8302.200000000001
$ws.Range("I3").Value = 1299.1428
$ws.Range("J3").Value = 24642.666
$ws.Range("K3").Value = 3897.4284
$ws.Range("L3").Value = 73927.99800000001
$ws.Range("M3").Value = -3785.4284
$ws.Range("N3").Value = -74151.99800000001

$ws.Range("H129").Value = 2152.8
$ws.Range("I129").Value = 647.4
$ws.Range("J129").Value = 2905.5
$ws.Range("K129").Value = 1942.2
$ws.Range("L129").Value = 8716.5
$ws.Range("M129").Value = 3057.8
$ws.Range("N129").Value = -18716.5

$ws.Range("H134").Value = 978.1429000000001
$ws.Range("I134").Value = 978.1429000000001
$ws.Range("K134").Value = 2934.4287
$ws.Range("M134").Value = 2135.5713

$ws.Range("H136").Value = 3115.85
$ws.Range("I136").Value = 2753.5789
$ws.Range("K136").Value = 8260.736699999999
$ws.Range("M136").Value = -3160.736699999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H57").Value = 25600
$ws.Range("J57").Value = 25600
$ws.Range("L57").Value = 25600
$ws.Range("N57").Value = -27240

$ws.Range("H126").Value = 836472.1
$ws.Range("I126").Value = 983620.25
$ws.Range("J126").Value = 2632.6667
$ws.Range("K126").Value = 2950860.75
$ws.Range("L126").Value = 7898.000100000001
$ws.Range("M126").Value = -2948390.75
$ws.Range("N126").Value = -12838.0001

$ws.Range("H132").Value = 4431679
$ws.Range("I132").Value = 5080846.5
$ws.Range("K132").Value = 15242539.5
$ws.Range("M132").Value = -15240009.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 9333
$ws.Range("I7").Value = 4000
$ws.Range("J7").Value = 19999
$ws.Range("K7").Value = 4000
$ws.Range("L7").Value = 19999
$ws.Range("M7").Value = -3888
$ws.Range("N7").Value = -20223

$ws.Range("H9").Value = 379.1111
$ws.Range("J9").Value = 619.75
$ws.Range("L9").Value = 619.75
$ws.Range("N9").Value = -1067.75

$ws.Range("H13").Value = 5926.5
$ws.Range("J13").Value = 6345.5
$ws.Range("L13").Value = 6345.5
$ws.Range("N13").Value = -6625.5

$ws.Range("H22").Value = 1024.6428
$ws.Range("I22").Value = 1024.6428
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 1024.6428
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = -729.6428000000001
$ws.Range("N22").ClearContents()

$ws.Range("H24").Value = 0
$ws.Range("I24").Value = 0
$ws.Range("J24").Value = 0
$ws.Range("K24").Value = 0
$ws.Range("L24").Value = 0
$ws.Range("N24").ClearContents()

$ws.Range("H27").Value = 1024.6428
$ws.Range("I27").Value = 1024.6428
$ws.Range("J27").Value = 0
$ws.Range("K27").Value = 1024.6428
$ws.Range("L27").Value = 0
$ws.Range("M27").Value = -917.6428000000001
$ws.Range("N27").ClearContents()

$ws.Range("H40").Value = 4527
$ws.Range("I40").Value = 3861.7693
$ws.Range("J40").Value = 5968.3335
$ws.Range("K40").Value = 3861.7693
$ws.Range("L40").Value = 5968.3335
$ws.Range("M40").Value = -3725.7693
$ws.Range("N40").Value = -6240.3335

$ws.Range("H100").Value = 22440.6
$ws.Range("J100").Value = 27300.75
$ws.Range("L100").Value = 27300.75
$ws.Range("N100").Value = -28382.75

$ws.Range("H122").Value = 3411.7646
$ws.Range("I122").Value = 3007.6924
$ws.Range("K122").Value = 9023.0772
$ws.Range("M122").Value = -6573.0772

$ws.Range("H126").Value = 9333
$ws.Range("I126").Value = 4000
$ws.Range("J126").Value = 19999
$ws.Range("K126").Value = 12000
$ws.Range("L126").Value = 59997
$ws.Range("M126").Value = -9530
$ws.Range("N126").Value = -64937

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H115").Value = 44377
$ws.Range("J115").Value = 44377
$ws.Range("L115").Value = 44377
$ws.Range("N115").Value = -47511

$ws.Range("H122").Value = 2510.8125
$ws.Range("I122").Value = 1817.8
$ws.Range("J122").Value = 3665.8333
$ws.Range("K122").Value = 5453.4
$ws.Range("L122").Value = 10997.4999
$ws.Range("M122").Value = -3003.4
$ws.Range("N122").Value = -15897.4999

$ws.Range("H126").Value = 3496.0715
$ws.Range("I126").Value = 3978.8333
$ws.Range("J126").Value = 599.5
$ws.Range("K126").Value = 11936.4999
$ws.Range("L126").Value = 1798.5
$ws.Range("M126").Value = -9466.499899999999
$ws.Range("N126").Value = -6738.5

$ws.Range("I132").Value = 7432507
$ws.Range("J132").Value = 47625720
$ws.Range("K132").Value = 22297521
$ws.Range("L132").Value = 142877160
$ws.Range("M132").Value = -22294991

$ws.Range("H136").Value = 17165804
$ws.Range("I136").Value = 21806860
$ws.Range("J136").Value = 29595.615
$ws.Range("K136").Value = 65420580
$ws.Range("L136").Value = 88786.845
$ws.Range("M136").Value = -65418030
$ws.Range("N136").Value = -93886.845
